# Generate Report for Handback
#
# This mirrors the localization tool's "handback" run: the zh-cn/de-de
# sheets gain a resolved "Latest Target File" + "Latest Handback File" +
# "Latest Handback DateTime" for each of the two rows (previously blank /
# "0001-01-01"), matching hyperlinks are added for the new target-file
# cells (same pattern as the existing source-file hyperlinks in column A),
# the "Ready for handoff" status flips to "Handed back: in sync with
# en-US" everywhere it appears, and several columns are widened so the
# new long file names / timestamps are readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 all shared this string)
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File / Latest Handback File /
#    Latest Handback DateTime for both rows, and hyperlink the target
#    file cell the same way the source file cell (column A) is linked.
# ---------------------------------------------------------------------
$zhcn.Range("I2").Value = "c799f5b3-20f1-42f0-8bf7-e29ed73f472f.md"
$zhcn.Range("J2").Value = "c799f5b3-20f1-42f0-8bf7-e29ed73f472f.7bfe6205cff3d1f893eecbdee7a5863609f0f013.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-19 18:38:35"

$zhcn.Range("I3").Value = "e768144c-e7db-41aa-8668-6e06049ad546.md"
$zhcn.Range("J3").Value = "e768144c-e7db-41aa-8668-6e06049ad546.f88d3ffaa0cd4eba99e11288d826f743ae1b704b.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-19 18:38:35"

$zhcn.Hyperlinks.Add(
    $zhcn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7b465c490a8746b83ee37214b9fda3bb172ca1d/e2e/c799f5b3-20f1-42f0-8bf7-e29ed73f472f.md",
    "",
    "",
    "c799f5b3-20f1-42f0-8bf7-e29ed73f472f.md"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7b465c490a8746b83ee37214b9fda3bb172ca1d/e2e/e768144c-e7db-41aa-8668-6e06049ad546.md",
    "",
    "",
    "e768144c-e7db-41aa-8668-6e06049ad546.md"
)

# ---------------------------------------------------------------------
# 3. de-de sheet: same shape of update, different (later) handback time.
# ---------------------------------------------------------------------
$dede.Range("I2").Value = "c799f5b3-20f1-42f0-8bf7-e29ed73f472f.md"
$dede.Range("J2").Value = "c799f5b3-20f1-42f0-8bf7-e29ed73f472f.7bfe6205cff3d1f893eecbdee7a5863609f0f013.de-de.xlf"
$dede.Range("K2").Value = "2016-08-19 18:38:41"

$dede.Range("I3").Value = "e768144c-e7db-41aa-8668-6e06049ad546.md"
$dede.Range("J3").Value = "e768144c-e7db-41aa-8668-6e06049ad546.f88d3ffaa0cd4eba99e11288d826f743ae1b704b.de-de.xlf"
$dede.Range("K3").Value = "2016-08-19 18:38:41"

$dede.Hyperlinks.Add(
    $dede.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7b465c490a8746b83ee37214b9fda3bb172ca1d/e2e/c799f5b3-20f1-42f0-8bf7-e29ed73f472f.md",
    "",
    "",
    "c799f5b3-20f1-42f0-8bf7-e29ed73f472f.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7b465c490a8746b83ee37214b9fda3bb172ca1d/e2e/e768144c-e7db-41aa-8668-6e06049ad546.md",
    "",
    "",
    "e768144c-e7db-41aa-8668-6e06049ad546.md"
)

# ---------------------------------------------------------------------
# 4. Column widths: widen the status/target/handback columns so the
#    new long values are fully visible.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527   # zh-cn status col
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527   # de-de status col

$zhcn.Columns.Item(3).ColumnWidth  = 29.9777047293527      # Status
$zhcn.Columns.Item(9).ColumnWidth  = 39.16666666666666     # Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = 39.16666666666666     # Latest Handback File

$dede.Columns.Item(3).ColumnWidth  = 29.9777047293527      # Status
$dede.Columns.Item(9).ColumnWidth  = 39.16666666666666     # Latest Target File
$dede.Columns.Item(10).ColumnWidth = 39.16666666666666     # Latest Handback File
